$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 2 - Table 1")

function Set-NumValue($cell, $val) {
    $cell.NumberFormat = "General"
    $cell.Value = $val
}

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}

# Row 3 - Jump
Set-NumValue $ws.Cells.Item(3,2) 0
Set-NumValue $ws.Cells.Item(3,3) 0
Set-NumValue $ws.Cells.Item(3,4) 0
Set-NumValue $ws.Cells.Item(3,5) 0
Set-NumValue $ws.Cells.Item(3,6) 0
Set-NumValue $ws.Cells.Item(3,7) 0
Set-NumValue $ws.Cells.Item(3,8) 0
Set-NumValue $ws.Cells.Item(3,9) 0
Set-NumValue $ws.Cells.Item(3,10) 1

# Row 4 - Branch
Set-NumValue $ws.Cells.Item(4,8) 0
Set-NumValue $ws.Cells.Item(4,9) 0
Set-NumValue $ws.Cells.Item(4,10) 2

# Row 5 - MemRead
Set-NumValue $ws.Cells.Item(5,6) 0
Set-NumValue $ws.Cells.Item(5,8) 0
Set-NumValue $ws.Cells.Item(5,9) 0
Set-NumValue $ws.Cells.Item(5,10) 2

# Row 6 - MemtoReg
Set-NumValue $ws.Cells.Item(6,6) 0
Set-NumValue $ws.Cells.Item(6,8) 0
Set-NumValue $ws.Cells.Item(6,9) 0
Set-NumValue $ws.Cells.Item(6,10) 2

# Row 7 - ALUOp
Set-TextValue $ws.Cells.Item(7,3) "000"
Set-TextValue $ws.Cells.Item(7,4) "000"
Set-TextValue $ws.Cells.Item(7,5) "000"
Set-NumValue $ws.Cells.Item(7,6) 101
Set-TextValue $ws.Cells.Item(7,7) "001"
Set-TextValue $ws.Cells.Item(7,8) "010"
Set-TextValue $ws.Cells.Item(7,9) "011"
Set-NumValue $ws.Cells.Item(7,10) 2

# Row 8 - MemWrite
Set-NumValue $ws.Cells.Item(8,8) 0
Set-NumValue $ws.Cells.Item(8,9) 0
Set-NumValue $ws.Cells.Item(8,10) 2

# Row 9 - ALUSrc
Set-NumValue $ws.Cells.Item(9,7) 1
Set-NumValue $ws.Cells.Item(9,8) 1
Set-NumValue $ws.Cells.Item(9,9) 1
Set-NumValue $ws.Cells.Item(9,10) 2

# Row 10 - RegWrite
Set-NumValue $ws.Cells.Item(10,8) 1
Set-NumValue $ws.Cells.Item(10,9) 1
Set-NumValue $ws.Cells.Item(10,10) 2
